$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 9095
$ws.Range("I46").Value = 5940
$ws.Range("J46").Value = 12250
$ws.Range("K46").Value = 17820
$ws.Range("L46").Value = 36750
$ws.Range("M46").Value = -17701
$ws.Range("N46").Value = -36988

$ws.Range("H60").Value = 9095
$ws.Range("I60").Value = 5940
$ws.Range("J60").Value = 12250
$ws.Range("K60").Value = 17820
$ws.Range("L60").Value = 36750
$ws.Range("M60").Value = -17336
$ws.Range("N60").Value = -37718

$ws.Range("H64").Value = 1039945.7
$ws.Range("I64").Value = 2178324.8
$ws.Range("J64").Value = 5055.636
$ws.Range("K64").Value = 2178324.8
$ws.Range("L64").Value = 5055.636
$ws.Range("M64").Value = -2178076.8
$ws.Range("N64").Value = -5551.636

$ws.Range("H67").Value = 1039945.7
$ws.Range("I67").Value = 2178324.8
$ws.Range("J67").Value = 5055.636
$ws.Range("K67").Value = 2178324.8
$ws.Range("L67").Value = 5055.636
$ws.Range("M67").Value = -2177466.8
$ws.Range("N67").Value = -6771.636

$ws.Range("H82").Value = 8038.375
$ws.Range("I82").Value = 4769
$ws.Range("K82").Value = 14307
$ws.Range("M82").Value = -13901

$ws.Range("H85").Value = 8038.375
$ws.Range("I85").Value = 4769
$ws.Range("K85").Value = 14307
$ws.Range("M85").Value = -12903

$ws.Range("H121").Value = 3729.1177
$ws.Range("J121").Value = 3887.5
$ws.Range("L121").Value = 11662.5
$ws.Range("N121").Value = -15156.5

$ws.Range("H124").Value = 58000
$ws.Range("J124").Value = 58000
$ws.Range("L124").Value = 58000
$ws.Range("N124").Value = -67820

$ws.Range("H132").Value = 7753.06
$ws.Range("I132").Value = 1969.1569
$ws.Range("J132").Value = 13773.041
$ws.Range("K132").Value = 5907.4707
$ws.Range("L132").Value = 41319.123
$ws.Range("M132").Value = -3377.4707
$ws.Range("N132").Value = -46379.123

$ws.Range("H137").Value = 3244.9565
$ws.Range("I137").Value = 2677.3333
$ws.Range("J137").Value = 4309.25
$ws.Range("K137").Value = 8031.999899999999
$ws.Range("L137").Value = 12927.75
$ws.Range("M137").Value = -5481.999899999999
$ws.Range("N137").Value = -18027.75

$ws.Range("H138").Value = 5511.4375
$ws.Range("I138").Value = 1345.8334
$ws.Range("J138").Value = 5942.3623
$ws.Range("K138").Value = 4037.5002
$ws.Range("L138").Value = 17827.0869
$ws.Range("M138").Value = 1102.4998
$ws.Range("N138").Value = -28107.0869

$ws.Range("H140").Value = 68462
$ws.Range("J140").Value = 68291.11
$ws.Range("L140").Value = 68291.11
$ws.Range("N140").Value = -78651.11

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5560.227
$ws.Range("I32").Value = 4442.8887
$ws.Range("K32").Value = 4442.8887
$ws.Range("M32").Value = -4155.8887

$ws.Range("H97").Value = 1326.6818
$ws.Range("I97").Value = 959.4
$ws.Range("K97").Value = 959.4
$ws.Range("M97").Value = -463.4

$ws.Range("H102").Value = 1055240
$ws.Range("I102").Value = 1055240
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1055240
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1053618
$ws.Range("N102").ClearContents()

$ws.Range("H132").Value = 37348.633
$ws.Range("I132").Value = 44059.855
$ws.Range("K132").Value = 132179.565
$ws.Range("M132").Value = -129649.565

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1987.2778
$ws.Range("I86").Value = 1635.6666
$ws.Range("J86").Value = 2338.889
$ws.Range("K86").Value = 1635.6666
$ws.Range("L86").Value = 2338.889
$ws.Range("M86").Value = -512.6666
$ws.Range("N86").Value = -4584.889

$ws.Range("H89").Value = 1987.2778
$ws.Range("I89").Value = 1635.6666
$ws.Range("J89").Value = 2338.889
$ws.Range("K89").Value = 8178.333000000001
$ws.Range("L89").Value = 11694.445
$ws.Range("M89").Value = -2562.333000000001
$ws.Range("N89").Value = -22926.445

$ws.Range("H94").Value = 1050722.6
$ws.Range("I94").Value = 1470721.8
$ws.Range("K94").Value = 1470721.8
$ws.Range("M94").Value = -1470270.8

$ws.Range("H107").Value = 27781072
$ws.Range("I107").Value = 27781072
$ws.Range("K107").Value = 27781072
$ws.Range("M107").Value = -27779152

$ws.Range("H134").Value = 5566.864
$ws.Range("I134").Value = 2832.7778
$ws.Range("K134").Value = 8498.3334
$ws.Range("M134").Value = -5963.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15154657
$ws.Range("I31").Value = 25001372
$ws.Range("K31").Value = 25001372
$ws.Range("M31").Value = -25001077

$ws.Range("H34").Value = 15154657
$ws.Range("I34").Value = 25001372
$ws.Range("K34").Value = 25001372
$ws.Range("M34").Value = -25001170

$ws.Range("H58").Value = 10000000
$ws.Range("I58").Value = 10000000
$ws.Range("K58").Value = 10000000
$ws.Range("M58").Value = -9999797

$ws.Range("H122").Value = 4810
$ws.Range("I122").Value = 4179.1763
$ws.Range("K122").Value = 12537.5289
$ws.Range("M122").Value = -10087.5289

$ws.Range("H132").Value = 21281144
$ws.Range("J132").Value = 21605.2
$ws.Range("L132").Value = 64815.60000000001
$ws.Range("N132").Value = -69875.60000000001

$ws.Range("H133").Value = 119826
$ws.Range("J133").Value = 119826
$ws.Range("L133").Value = 119826
$ws.Range("N133").Value = -124886

$ws.Range("H136").Value = 10000000
$ws.Range("I136").Value = 10000000
$ws.Range("K136").Value = 30000000
$ws.Range("M136").Value = -29997450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4140.7144
$ws.Range("I32").Value = 4058.4
$ws.Range("J32").Value = 4346.5
$ws.Range("K32").Value = 12175.2
$ws.Range("L32").Value = 13039.5
$ws.Range("M32").Value = -11892.2
$ws.Range("N32").Value = -13605.5

$ws.Range("H62").Value = 12402.6
$ws.Range("J62").Value = 12402.6
$ws.Range("L62").Value = 37207.8
$ws.Range("N62").Value = -38579.8

$ws.Range("H65").Value = 12402.6
$ws.Range("J65").Value = 12402.6
$ws.Range("L65").Value = 111623.4
$ws.Range("N65").Value = -118487.4

$ws.Range("H128").Value = 343293.66
$ws.Range("I128").Value = 343293.66
$ws.Range("K128").Value = 1029880.98
$ws.Range("M128").Value = -1024900.98

$ws.Range("H131").Value = 14356116
$ws.Range("J131").Value = 6631118
$ws.Range("L131").Value = 19893354
$ws.Range("N131").Value = -19903434

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H70").Value = 1648294.8
$ws.Range("I70").Value = 3180039.8
$ws.Range("J70").Value = 7139.5
$ws.Range("K70").Value = 3180039.8
$ws.Range("L70").Value = 7139.5
$ws.Range("M70").Value = -3179769.8
$ws.Range("N70").Value = -7679.5

$ws.Range("H73").Value = 1648294.8
$ws.Range("I73").Value = 3180039.8
$ws.Range("J73").Value = 7139.5
$ws.Range("K73").Value = 3180039.8
$ws.Range("L73").Value = 7139.5
$ws.Range("M73").Value = -3179103.8
$ws.Range("N73").Value = -9011.5

$ws.Range("H80").Value = 1363518.5
$ws.Range("I80").Value = 1963126.9
$ws.Range("K80").Value = 1963126.9
$ws.Range("M80").Value = -1962128.9

$ws.Range("H83").Value = 1363518.5
$ws.Range("I83").Value = 1963126.9
$ws.Range("K83").Value = 9815634.5
$ws.Range("M83").Value = -9810642.5

$ws.Range("H102").Value = 14289816
$ws.Range("I102").Value = 18520914
$ws.Range("J102").Value = 9860.125
$ws.Range("K102").Value = 18520914
$ws.Range("L102").Value = 9860.125
$ws.Range("M102").Value = -18519292
$ws.Range("N102").Value = -13104.125

$ws.Range("H126").Value = 3548.1072
$ws.Range("I126").Value = 2102.1304
$ws.Range("K126").Value = 6306.3912
$ws.Range("M126").Value = -3836.3912

$ws.Range("H132").Value = 3164.362
$ws.Range("I132").Value = 3149.7778
$ws.Range("J132").Value = 3214.8462
$ws.Range("K132").Value = 9449.3334
$ws.Range("L132").Value = 9644.5386
$ws.Range("M132").Value = -6919.3334
$ws.Range("N132").Value = -14704.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1750057.8
$ws.Range("I68").Value = 3789505.2
$ws.Range("K68").Value = 3789505.2
$ws.Range("M68").Value = -3788756.2

$ws.Range("H71").Value = 1750057.8
$ws.Range("I71").Value = 3789505.2
$ws.Range("K71").Value = 18947526
$ws.Range("M71").Value = -18943782

$ws.Range("H122").Value = 7489.517
$ws.Range("I122").Value = 4254.278
$ws.Range("J122").Value = 12783.546
$ws.Range("K122").Value = 12762.834
$ws.Range("L122").Value = 38350.638
$ws.Range("M122").Value = -10312.834
$ws.Range("N122").Value = -43250.638

$ws.Range("H132").Value = 4238.7144
$ws.Range("I132").Value = 4238.7144
$ws.Range("K132").Value = 12716.1432
$ws.Range("M132").Value = -10186.1432

$ws.Range("H136").Value = 3988.17
$ws.Range("I136").Value = 3987.949
$ws.Range("J136").Value = 3999
$ws.Range("K136").Value = 11963.847
$ws.Range("L136").Value = 11997
$ws.Range("M136").Value = -9413.847
$ws.Range("N136").Value = -17097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3483915.8
$ws.Range("I132").Value = 6962583.5
$ws.Range("K132").Value = 20887750.5
$ws.Range("M132").Value = -20885220.5

$ws.Range("H136").Value = 7759.87
$ws.Range("I136").Value = 3578.25
$ws.Range("J136").Value = 9080.382
$ws.Range("K136").Value = 10734.75
$ws.Range("L136").Value = 27241.146
$ws.Range("M136").Value = -8184.75
$ws.Range("N136").Value = -32341.146

$ws.Range("H139").Value = 98085.39999999999
$ws.Range("J139").Value = 98085.39999999999
$ws.Range("L139").Value = 98085.39999999999
$ws.Range("N139").Value = -108365.4
